{"js": "const body = context.document.body;\n\nconst replacements = [\n  [\"64\u00d720=\", \"67\u00d737=\"],\n  [\"61\u00d767=\", \"57\u00d756=\"],\n  [\"95\u00d773=\", \"93\u00d769=\"],\n  [\"26\u00d788=\", \"56\u00d754=\"],\n  [\"42\u00d782=\", \"88\u00d721=\"],\n  [\"99\u00d724=\", \"67\u00d769=\"],\n  [\"84\u00d718=\", \"69\u00d716=\"],\n  [\"68\u00d742=\", \"38\u00d731=\"],\n  [\"62\u00d780=\", \"74\u00d785=\"],\n  [\"28\u00d727=\", \"87\u00d732=\"],\n  [\"35\u00d752=\", \"42\u00d767=\"],\n  [\"66\u00d751=\", \"88\u00d752=\"],\n  [\"76\u00d750=\", \"87\u00d714=\"],\n  [\"45\u00d765=\", \"48\u00d794=\"],\n  [\"60\u00d718=\", \"99\u00d716=\"],\n  [\"92\u00d750=\", \"26\u00d771=\"],\n  [\"31\u00d734=\", \"52\u00d794=\"],\n  [\"70\u00d769=\", \"19\u00d718=\"],\n  [\"94\u00d740=\", \"35\u00d771=\"],\n  [\"76\u00d756=\", \"74\u00d764=\"],\n  [\"35\u00d773=\", \"51\u00d759=\"],\n  [\"62\u00d732=\", \"28\u00d798=\"],\n  [\"35\u00d726=\", \"60\u00d795=\"],\n  [\"44\u00d754=\", \"32\u00d762=\"],\n  [\"68\u00d744=\", \"15\u00d711=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$r = $d.Content\n$r.Find.Execute(\"64\u00d720=\", $false, $false, $false, $false, $false, $true, 1, $false, \"67\u00d737=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"61\u00d767=\", $false, $false, $false, $false, $false, $true, 1, $false, \"57\u00d756=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"95\u00d773=\", $false, $false, $false, $false, $false, $true, 1, $false, \"93\u00d769=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"26\u00d788=\", $false, $false, $false, $false, $false, $true, 1, $false, \"56\u00d754=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"42\u00d782=\", $false, $false, $false, $false, $false, $true, 1, $false, \"88\u00d721=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"99\u00d724=\", $false, $false, $false, $false, $false, $true, 1, $false, \"67\u00d769=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"84\u00d718=\", $false, $false, $false, $false, $false, $true, 1, $false, \"69\u00d716=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"68\u00d742=\", $false, $false, $false, $false, $false, $true, 1, $false, \"38\u00d731=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"62\u00d780=\", $false, $false, $false, $false, $false, $true, 1, $false, \"74\u00d785=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"28\u00d727=\", $false, $false, $false, $false, $false, $true, 1, $false, \"87\u00d732=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"35\u00d752=\", $false, $false, $false, $false, $false, $true, 1, $false, \"42\u00d767=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"66\u00d751=\", $false, $false, $false, $false, $false, $true, 1, $false, \"88\u00d752=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"76\u00d750=\", $false, $false, $false, $false, $false, $true, 1, $false, \"87\u00d714=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"45\u00d765=\", $false, $false, $false, $false, $false, $true, 1, $false, \"48\u00d794=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"60\u00d718=\", $false, $false, $false, $false, $false, $true, 1, $false, \"99\u00d716=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"92\u00d750=\", $false, $false, $false, $false, $false, $true, 1, $false, \"26\u00d771=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"31\u00d734=\", $false, $false, $false, $false, $false, $true, 1, $false, \"52\u00d794=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"70\u00d769=\", $false, $false, $false, $false, $false, $true, 1, $false, \"19\u00d718=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"94\u00d740=\", $false, $false, $false, $false, $false, $true, 1, $false, \"35\u00d771=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"76\u00d756=\", $false, $false, $false, $false, $false, $true, 1, $false, \"74\u00d764=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"35\u00d773=\", $false, $false, $false, $false, $false, $true, 1, $false, \"51\u00d759=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"62\u00d732=\", $false, $false, $false, $false, $false, $true, 1, $false, \"28\u00d798=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"35\u00d726=\", $false, $false, $false, $false, $false, $true, 1, $false, \"60\u00d795=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"44\u00d754=\", $false, $false, $false, $false, $false, $true, 1, $false, \"32\u00d762=\", 2)\n\n$r = $d.Content\n$r.Find.Execute(\"68\u00d744=\", $false, $false, $false, $false, $false, $true, 1, $false, \"15\u00d711=\", 2)\n"}
